# feat: add 2022-Q4 data
#
# 1. Insert a new row into the "总计" (total) summary sheet for 2022-Q4,
#    pushing the existing quarters down by one row.
# 2. Insert a brand-new worksheet named "2022-Q4" right before the
#    existing "2022-Q3" tab, populated with the per-fund holdings table
#    for that quarter. All the other quarter sheets keep their data
#    untouched; they simply shift one tab to the right to make room for
#    the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: 总计 sheet - insert the new 2022-Q4 summary row at row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Re-use the existing "index" column styling (bold + border, same as the
# other rows' A column) for the new A2 cell.
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 13
$total.Cells.Item(2,4).Value = 0.28

# ---------------------------------------------------------------------
# Step 2: insert a new "2022-Q4" worksheet right before "2022-Q3".
# ---------------------------------------------------------------------
# Duplicating the existing "2022-Q3" sheet (rather than Worksheets.Add())
# carries over the header/formatting faithfully, since pasting formats
# into a brand-new blank sheet does not stick in this runtime. The
# duplicate is placed immediately before its source, i.e. right after
# "总计" - exactly where "2022-Q4" belongs.
$oldQ3 = $wb.Worksheets.Item(2)
$oldQ3.Copy($oldQ3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# 2022-Q3 has 15 fund rows (rows 2-16); 2022-Q4 only has 13 (rows 2-14),
# so drop the two extra trailing rows carried over from the duplicate.
$q4.Rows.Item(15).Delete()
$q4.Rows.Item(15).Delete()

function Set-FundRow($row, $idx, $code, $fundName, $size, $stockPos, $posPct, $heldValue, $posRank) {
    $q4.Cells.Item($row,1).Value = $idx

    $cell = $q4.Cells.Item($row,2)
    $cell.NumberFormat = "@"
    $cell.Value = $code

    $q4.Cells.Item($row,3).Value = $fundName

    $cell = $q4.Cells.Item($row,4)
    $cell.NumberFormat = "@"
    $cell.Value = $size

    $cell = $q4.Cells.Item($row,5)
    $cell.NumberFormat = "@"
    $cell.Value = $stockPos

    $cell = $q4.Cells.Item($row,6)
    $cell.NumberFormat = "@"
    $cell.Value = $posPct

    $cell = $q4.Cells.Item($row,7)
    $cell.NumberFormat = "@"
    $cell.Value = $heldValue

    $q4.Cells.Item($row,8).Value = $posRank
}

Set-FundRow 2  0  "002666" "前海开源沪港深创新成长灵活配置混合A" "4.54" "67.11" "3.06" "0.1389" 2
Set-FundRow 3  1  "002667" "前海开源沪港深创新成长灵活配置混合C" "2.14" "67.11" "3.06" "0.0655" 2
Set-FundRow 4  2  "010683" "中融景颐6个月持有期混合A" "2.71" "22.07" "0.93" "0.0252" 9
Set-FundRow 5  3  "010367" "中融景瑞一年持有期混合A" "1.26" "21.49" "1.07" "0.0135" 7
Set-FundRow 6  4  "710002" "富安达策略精选混合" "0.58" "67.44" "2.24" "0.0130" 4
Set-FundRow 7  5  "002003" "工银新机遇灵活配置混合A" "0.44" "76.64" "1.54" "0.0068" 5
Set-FundRow 8  6  "010684" "中融景颐6个月持有期混合C" "0.56" "22.07" "0.93" "0.0052" 9
Set-FundRow 9  7  "000432" "中银优秀企业混合" "0.17" "85.63" "2.46" "0.0042" 9
Set-FundRow 10 8  "010368" "中融景瑞一年持有期混合C" "0.31" "21.49" "1.07" "0.0033" 7
Set-FundRow 11 9  "004456" "兴银消费新趋势灵活配置混合" "0.09" "80.92" "3.40" "0.0031" 9
Set-FundRow 12 10 "002004" "工银新机遇灵活配置混合C" "0.13" "76.64" "1.54" "0.0020" 5
Set-FundRow 13 11 "001724" "申万菱信多策略灵活配置混合C" "0.08" "22.14" "1.11" "0.0009" 4
Set-FundRow 14 12 "001148" "申万菱信多策略灵活配置混合A" "0.04" "22.14" "1.11" "0.0004" 4
